$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.792.61"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.280.08"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.645"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.24%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.645"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0977"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "2.622.86"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "2.278.87"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "42.678.34"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "236.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "167.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0872"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.68"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  +2.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  -4.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.44%  "
$ws.Range("E47").Value = "  -1.02%  "
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -1.32%  "
$ws.Range("E51").Value = "  -1.31%  "

# Row 43/44 swap: FraxShare <-> MultiversX
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "61.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
